$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking values are stored as literal text (matching the
# original inlineStr cell contents) rather than being parsed into numbers.
$textCells = @("D2", "D3", "D4", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D27", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "243.01"
$ws.Range("D3").Value = "23.09"
$ws.Range("D4").Value = "5.429"
$ws.Range("D6").Value = "3.447"
$ws.Range("D7").Value = "6.538"
$ws.Range("D8").Value = "0.8114"
$ws.Range("D9").Value = "0.9737"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "0.01130"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1418"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.07452"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "0.03259"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "0.03063"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "0.09344"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "3.873"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "0.001575"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "0.04674"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("D19").Value = "0.005878"
$ws.Range("D20").Value = "0.001256"
$ws.Range("D21").Value = "0.004903"
$ws.Range("D22").Value = "0.00006811"
$ws.Range("D23").Value = "3.591"
$ws.Range("D24").Value = "2.133"
$ws.Range("D27").Value = "0.0002288"
$ws.Range("D40").Value = "0.03936"
$ws.Range("D41").Value = "0.006188"
$ws.Range("D42").Value = "0.1071"
$ws.Range("D43").Value = "0.003005"
$ws.Range("D44").Value = "0.009157"
$ws.Range("E44").Value = "43LocalTradersLCT"
$ws.Range("D45").Value = "0.00005186"
$ws.Range("D47").Value = "0.7512"
$ws.Range("D48").Value = "0.002294"
$ws.Range("D50").Value = "0.0002003"
